$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Updated values (row-by-row) from prepare & render with final data
$ws.Range("B2").Value = 0.556007222541788
$ws.Range("C2").Value = 0.612606004275779
$ws.Range("K2").Value = 0.438282499020205
$ws.Range("L2").Value = 0.595756191953926
$ws.Range("N2").Value = 0.508231644030169

$ws.Range("B3").Value = 0.50312291438834
$ws.Range("K3").Value = 0.351050882580874
$ws.Range("L3").Value = 0.604399737467109
$ws.Range("N3").Value = 0.447300643788012

$ws.Range("B4").Value = 0.680881448179833
$ws.Range("K4").Value = 0.758076861129753
$ws.Range("L4").Value = 0.688615273248795
$ws.Range("N4").Value = 0.616918649447641

$ws.Range("B5").Value = 0.609601586795904
$ws.Range("K5").Value = 0.571096670838126
$ws.Range("L5").Value = 0.778963825426238
$ws.Range("N5").Value = 0.472900191628792

$ws.Range("B6").Value = 0.675595447215337
$ws.Range("K6").Value = 0.557841849059486
$ws.Range("N6").Value = 0.666651932459956

$ws.Range("B7").Value = 0.364717906507653
$ws.Range("K7").Value = 0.22156020948145
$ws.Range("N7").Value = 0.372217577193357

$ws.Range("B8").Value = 0.347853243460036
$ws.Range("N8").Value = 0.306445646731996

$ws.Range("B9").Value = 0.410626908494325
$ws.Range("K9").Value = 0.319383802321488
$ws.Range("L9").Value = 0.389233362357354
$ws.Range("N9").Value = 0.40055514051731
